$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture existing data rows (2-21, columns A-C) before shifting
$oldValues = @()
for ($r = 2; $r -le 21; $r++) {
    $rowVals = @($ws.Cells.Item($r, 1).Value2, $ws.Cells.Item($r, 2).Value2, $ws.Cells.Item($r, 3).Value2)
    $oldValues += ,$rowVals
}

# Write old data back shifted down by 6 rows (new rows 8-27), without using Insert
# so no formatting gets copied/propagated.
for ($i = 0; $i -lt $oldValues.Length; $i++) {
    $row = 8 + $i
    $ws.Cells.Item($row, 1).Value = $oldValues[$i][0]
    $ws.Cells.Item($row, 2).Value = $oldValues[$i][1]
    $ws.Cells.Item($row, 3).Value = $oldValues[$i][2]
}

# New rows inserted at the top (rows 2-7)
$newTop = @(
    @(-0.0386372283101081, 0.008399397134780801, -0.0021380283869802),
    @(0.0134390350431203, 0.0704022198915481, -0.0390953756868839),
    @(0.0154243474826216, 0.030695978552103, 0.0099265603348612),
    @(0.041233405470848, -0.0022907445672899, 0.0502436682581901),
    @(-0.030695978552103, -0.062460970133543, 0.0204639863222837),
    @(0.0216857157647609, -0.0343611687421798, -0.0035124751739203)
)

for ($i = 0; $i -lt $newTop.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $newTop[$i][0]
    $ws.Cells.Item($row, 2).Value = $newTop[$i][1]
    $ws.Cells.Item($row, 3).Value = $newTop[$i][2]
}

# New rows appended at the bottom (rows 28-31)
$newBottom = @(
    @(0.0047342055477201, 0.0751364231109619, 0.0545197241008281),
    @(0.0003054326225537, 0.0100792767480015, -0.0006108652451075),
    @(-0.0152716310694813, -0.00167987938039, 0.047036625444889),
    @(-0.0091629782691597, -0.0114537235349416, 0.0062613687478005)
)

for ($i = 0; $i -lt $newBottom.Length; $i++) {
    $row = 28 + $i
    $ws.Cells.Item($row, 1).Value = $newBottom[$i][0]
    $ws.Cells.Item($row, 2).Value = $newBottom[$i][1]
    $ws.Cells.Item($row, 3).Value = $newBottom[$i][2]
}
